$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.120.73'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '3.391.28'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.388.13'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.26%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('D13').Value = '3.967.68'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').Value = '3.387.95'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '61.152.22'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.551'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.511.03'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.26%  '
$ws.Range('E29').Value = '  +8.56%  '
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.56'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.13'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '165.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('E47').Value = '  -3.34%  '
$ws.Range('D48').Value = '2.526.73'
$ws.Range('E48').Value = '  +7.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.78'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.71%  '
